$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
